$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while forcing text storage (no numeric/date
# auto-conversion) and without leaving the cells visible number format
# or style changed - we snapshot + restore Style around the write.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $savedStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $savedStyle
}

# Row 2
Set-TextValue 'D2' '29.685.97'
Set-TextValue 'E2' '  -1.37%  '
# Row 3
Set-TextValue 'D3' '2.096.36'
Set-TextValue 'E3' '  -0.51%  '
# Row 4
Set-TextValue 'E4' '  +0.55%  '
# Row 5
Set-TextValue 'D5' '343.27'
Set-TextValue 'E5' '  -1.90%  '
# Row 7
Set-TextValue 'D7' '0.5168'
Set-TextValue 'E7' '  +0.11%  '
# Row 8
Set-TextValue 'D8' '0.4373'
Set-TextValue 'E8' '  -2.37%  '
# Row 9
Set-TextValue 'D9' '53.61'
Set-TextValue 'E9' '  +1.52%  '
# Row 10
Set-TextValue 'D10' '0.09167'
Set-TextValue 'E10' '  +2.31%  '
# Row 11
Set-TextValue 'E11' '  -0.99%  '
# Row 12
Set-TextValue 'D12' '24.57'
Set-TextValue 'E12' '  -5.30%  '
# Row 13
Set-TextValue 'D13' '2.067.80'
Set-TextValue 'E13' '  -1.78%  '
# Row 14
Set-TextValue 'D14' '6.757'
# Row 15
Set-TextValue 'D15' '8.138'
Set-TextValue 'E15' '  -1.55%  '
# Row 16
Set-TextValue 'D16' '102.33'
Set-TextValue 'E16' '  +3.06%  '
# Row 17
Set-TextValue 'D17' '0.00001150'
Set-TextValue 'E17' '  +0.11%  '
# Row 18
Set-TextValue 'D18' '1.010'
Set-TextValue 'E18' '  +0.52%  '
# Row 19
Set-TextValue 'D19' '21.01'
Set-TextValue 'E19' '  +0.54%  '
# Row 20
Set-TextValue 'D20' '0.06671'
Set-TextValue 'E20' '  +0.00%  '
# Row 21
Set-TextValue 'D21' '1.008'
Set-TextValue 'E21' '  +0.52%  '
# Row 22
Set-TextValue 'D22' '6.197'
Set-TextValue 'E22' '  -1.64%  '
# Row 23
Set-TextValue 'D23' '29.750.41'
Set-TextValue 'E23' '  -1.48%  '
# Row 24
Set-TextValue 'D24' '12.66'
Set-TextValue 'E24' '  -1.83%  '
# Row 25
Set-TextValue 'D25' '2.305'
Set-TextValue 'E25' '  -2.05%  '
# Row 26
Set-TextValue 'D26' '2.327.50'
Set-TextValue 'E26' '  -1.28%  '
# Row 27
Set-TextValue 'D27' '21.87'
Set-TextValue 'E27' '  -1.05%  '
# Row 28
Set-TextValue 'D28' '161.87'
Set-TextValue 'E28' '  -0.73%  '
# Row 29
Set-TextValue 'D29' '2.490'
Set-TextValue 'E29' '  -3.05%  '
# Row 30
Set-TextValue 'D30' '133.20'
Set-TextValue 'E30' '  -0.43%  '
# Row 31
Set-TextValue 'D31' '1.127'
Set-TextValue 'E31' '  -4.85%  '
# Row 32
Set-TextValue 'D32' '1.667'
Set-TextValue 'E32' '  +1.22%  '
# Row 33
Set-TextValue 'E33' '  -1.87%  '
# Row 34
Set-TextValue 'D34' '6.183'
Set-TextValue 'E34' '  -1.46%  '
# Row 35
Set-TextValue 'D35' '3.964'
Set-TextValue 'E35' '  -0.53%  '
# Row 36
Set-TextValue 'D36' '6.301'
Set-TextValue 'E36' '  +6.49%  '
# Row 37
Set-TextValue 'D37' '10.43'
Set-TextValue 'E37' '  +1.87%  '
# Row 38
Set-TextValue 'D38' '0.02577'
Set-TextValue 'E38' '  -0.48%  '
# Row 39
Set-TextValue 'D39' '0.06692'
Set-TextValue 'E39' '  -2.23%  '
# Row 40
Set-TextValue 'D40' '0.6986'
Set-TextValue 'E40' '  +1.86%  '
# Row 41
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D41' '1.339'
Set-TextValue 'E41' '  +6.92%  '
# Row 42
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D42' '12.42'
Set-TextValue 'E42' '  -1.40%  '
# Row 43
Set-TextValue 'D43' '0.2210'
Set-TextValue 'E43' '  -4.91%  '
# Row 44
Set-TextValue 'D44' '0.6782'
Set-TextValue 'E44' '  +5.40%  '
# Row 45
Set-TextValue 'D45' '14.36'
Set-TextValue 'E45' '  +0.62%  '
# Row 46
Set-TextValue 'D46' '2.314'
Set-TextValue 'E46' '  +0.01%  '
# Row 47
Set-TextValue 'D47' '0.00000000362'
Set-TextValue 'E47' '  -1.64%  '
# Row 48
Set-TextValue 'D48' '3.619'
Set-TextValue 'E48' '  -1.28%  '
# Row 49
Set-TextValue 'D49' '1.207'
Set-TextValue 'E49' '  +3.33%  '
# Row 50
Set-TextValue 'E50' '  -0.83%  '
# Row 51
Set-TextValue 'D51' '81.02'
Set-TextValue 'E51' '  -3.59%  '
